# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting already used by the other header cells (e.g. G1 "sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's formatting (bold/border/alignment style) onto H1 so the new
# header cell matches the existing ones exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header + data values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
